# Apply "added phyre as test" edit to the tool_list workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the tags text for FAIRDOMHub (row 13, column E): add a space after the
#    first comma so it reads "collect, process, analyse, share, preserve, storage, privacy".
$ws.Range("E13").Value = "collect, process, analyse, share, preserve, storage, privacy"

# 2. Append a new row (16) for the Phyre2 tool.
$ws.Range("A16").Value = "Phyre2"
$ws.Range("B16").Value = "http://www.sbg.bio.ic.ac.uk/~phyre2"
$ws.Range("C16").Value = "Protein Homology/analogY Recognition Engine"
$ws.Range("D16").Value = "phyre"
$ws.Range("E16").Value = "process, analyse"

# Register the hyperlink, then apply the same "Hyperlink" style used by the
# other link cells (must be set after Hyperlinks.Add so it sticks).
$ws.Hyperlinks.Add($ws.Range("B16"), "http://www.sbg.bio.ic.ac.uk/~phyre2") | Out-Null
$ws.Range("B16").Style = "Hyperlink"

# Match the author's final selection (cell E18) as recorded in the saved file.
$ws.Range("E18").Select() | Out-Null
